# Restored from revision of admin on 06/03/2020 07:23:47 AM.TEST
# Author: admin. Type: SAVE.
#
# Semantic change: cell C10 on the "Rules" sheet changes from 18 to 1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C10").Value = 1
